$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Clear()

$ws.Range("A1").Value = "rishabh"

$ws.Range("G6").Select()
